$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean: wipe all cell contents (keeps existing formatting/styles in place)
$ws.Cells.ClearContents()

# Re-write header row (row 1) - unchanged text, same left-to-right order
$headers = @("Department","Subcategory","Category","Article Number","CodingType","UOMName","HSNCode","ExtDescription","Description","Color","Style","Size","Brand","Supplier","ItemCode","ItemId","Pur Price","ItemMRP","ItemWSP","Quantity","InvoiceNo","InvoiceDt")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Write the two data rows column-by-column (column-major) so newly introduced
# shared strings are appended in the same order as the target workbook.

# Column A - Department
$ws.Range("A2").Value = "TROUSER"
$ws.Range("A3").Value = "SHIRT"

# Column B - Subcategory
$ws.Range("B2").Value = "S3"
$ws.Range("B3").Value = "S2"

# Column C - Category
$ws.Range("C2").Value = "C2"
$ws.Range("C3").Value = "C4"

# Column D - Article Number
$ws.Range("D2").Value = "Sidhhit4114TSRTSR001C1S3"
$ws.Range("D3").Value = "suffix3612suffixST34C1S2"

# Column E - CodingType (numeric)
$ws.Range("E2").Value = 3
$ws.Range("E3").Value = 3

# Column F - UOMName
$ws.Range("F2").Value = "pcs"
$ws.Range("F3").Value = "pcs"

# Column G - HSNCode (row 3 has no HSNCode)
$ws.Range("G2").Value = "TSR1200987ZZ"

# Column I - Description (numeric)
$ws.Range("I2").Value = 6
$ws.Range("I3").Value = 7

# Column J - Color
$ws.Range("J2").Value = "L GREY"
$ws.Range("J3").Value = "OUTFIT BLACK"

# Column L - Size
$ws.Range("L2").Value = "L"
$ws.Range("L3").Value = "XL"

# Column M - Brand (numeric-looking text, force text storage)
$ws.Range("M2").Value = "'123"
$ws.Range("M2").Style = "Normal"
$ws.Range("M3").Value = "'1111"
$ws.Range("M3").Style = "Normal"

# Column N - Supplier
$ws.Range("N2").Value = "Siddhivinayak Apparel"
$ws.Range("N3").Value = "XYZ"

# Column Q - Pur Price (numeric-looking text, force text storage)
$ws.Range("Q2").Value = "'1000"
$ws.Range("Q2").Style = "Normal"
$ws.Range("Q3").Value = "'2000"
$ws.Range("Q3").Style = "Normal"

# Column R - ItemMRP (numeric-looking text, force text storage)
$ws.Range("R2").Value = "'2195.00"
$ws.Range("R2").Style = "Normal"
$ws.Range("R3").Value = "'111"
$ws.Range("R3").Style = "Normal"

# Column T - Quantity (numeric)
$ws.Range("T2").Value = 2
$ws.Range("T3").Value = 2

# Column U - InvoiceNo
$ws.Range("U2").Value = "INVOICE1221"
$ws.Range("U3").Value = "INVOICE4533"

# Column V - InvoiceDt (date serials, keeps existing date style)
$ws.Range("V2").Value = 45433
$ws.Range("V3").Value = 45421

# Remove row 4 entirely (workbook shrank from 3 data rows to 2)
$ws.Rows(4).Delete()
